# fix bug import + search + file data
#
# The "Avatar" column (J) on each sheet held relative paths like
# "Images/<brand>/avatarNN.jpg" pointing at a local image-import folder.
# The import/search code actually needs an absolute path, so every avatar
# cell is rewritten to "D:\Images\<brand>\avatarNN.jpg" (same NN, just a
# different path shape). Also nudges each sheet's selection/active-tab
# state to reflect where the user was poking around while testing the fix
# (ending up on the Xiaomi tab).

$wb = $excel.ActiveWorkbook

function Set-AvatarColumn {
    param($SheetName, $Brand, $RowOrder)

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in $RowOrder) {
        $n = (($row - 4) % 10) + 1
        $fileName = "avatar{0:D2}.jpg" -f $n
        $ws.Cells.Item($row, 10).Value = "D:\Images\$Brand\$fileName"
    }
}

# Row-write order reproduces the same cell-entry sequence the original
# author used (rows aren't touched strictly top-to-bottom on every sheet).
Set-AvatarColumn "iPhone"  "iphone"  @(5,6,7,8,9,10,11,12,13,4,14,15,16,17)
Set-AvatarColumn "Samsung" "samsung" @(5,4,6,7,8,9,10,11,12,13,14,15,16,17)
Set-AvatarColumn "Xiaomi"  "xiaomi"  @(4,5,6,7,8,9,10,11,12,13,14,15,16,17)

# Replay the selection trail across the three sheets, ending on Xiaomi
# (matches the saved workbook view / activeTab state).
$wsIPhone = $wb.Worksheets.Item("iPhone")
$wsIPhone.Select()
$wsIPhone.Range("J4").Select()

$wsSamsung = $wb.Worksheets.Item("Samsung")
$wsSamsung.Select()
$wsSamsung.Range("J6").Select()

$wsXiaomi = $wb.Worksheets.Item("Xiaomi")
$wsXiaomi.Select()
$wsXiaomi.Range("J5").Select()
